# Update Name of Algo
# Apply updated imputation result values to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E7").Value = 13.417
$ws.Range("A8").Value = -21.107
$ws.Range("A10").Value = -20.945
$ws.Range("A12").Value = -21.452
$ws.Range("C13").Value = -12.686
$ws.Range("A18").Value = -21.766
$ws.Range("E20").Value = 12.932
$ws.Range("A25").Value = -21.58
